# "Started writting report, need to test perfect branch one more time"
#
# 1) C58 held the placeholder shared-string " 0.012728765345678653" (the
#    "perfect branch" value jotted down as text); replace it with the real
#    numeric result so the TEST 4 average recalculates correctly.
# 2) Append a new "TEST 9 / ORIGINAL CW" results block (rows 145-160),
#    mirroring the layout used by every other test block on the sheet.
# 3) Update the view state to where the report-writer had scrolled to.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- 1. TEST 4 / C58: turn the leftover text placeholder into a real number
$ws.Range("C58").Value = 0.0127287653456786

# --- 2. New block: TEST 9 / ORIGINAL CW -----------------------------------
$ws.Range("A145").Value = "TEST 9"
$ws.Range("B145").Value = "ORIGINAL CW"

$ws.Range("A147").Value = "Method:"
$ws.Range("B147").Value = "5 hidden nodes, 200 population, "
$ws.Range("B148").Value = "replace worst, random parent selection,"

$ws.Range("A149").Value = "Test No"
$ws.Range("B149").Value = "Training Fitness"
$ws.Range("C149").Value = "Test Fitness"

$ws.Range("A150").Value = 1
$ws.Range("B150").Value = 0.13299845414636799
$ws.Range("C150").Value = 0.27902066547717003

$ws.Range("A151").Value = 2
$ws.Range("B151").Value = 0.13563810204085
$ws.Range("C151").Value = 0.28723236699432197

$ws.Range("A152").Value = 3
$ws.Range("B152").Value = 0.120337625475382
$ws.Range("C152").Value = 0.27564133468400798

$ws.Range("A153").Value = 4
$ws.Range("B153").Value = 0.17801790634224299
$ws.Range("C153").Value = 0.35197137499534698

$ws.Range("A154").Value = 5
$ws.Range("B154").Value = 0.12931036207077801
$ws.Range("C154").Value = 0.23684812577423101

$ws.Range("A155").Value = 6
$ws.Range("B155").Value = 0.20480812931675599
$ws.Range("C155").Value = 0.35089827423202502

$ws.Range("A156").Value = 7
$ws.Range("B156").Value = 0.090942183983095698
$ws.Range("C156").Value = 0.21102485324325099

$ws.Range("A157").Value = 8
$ws.Range("B157").Value = 0.092871148175299095
$ws.Range("C157").Value = 0.15266968401499001

$ws.Range("A158").Value = 9
$ws.Range("B158").Value = 0.116326957866423
$ws.Range("C158").Value = 0.29704662278780902

$ws.Range("A159").Value = 10
$ws.Range("B159").Value = 0.174028696284241
$ws.Range("C159").Value = 0.29677895377643898

$ws.Range("A160").Value = "Average:"
$ws.Range("B160").Formula = "=AVERAGE(B150:B159)"
$ws.Range("C160").Formula = "=AVERAGE(C150:C159)"

# --- 3. View state: scrolled down to the new block, M59 selected ----------
$ws.Range("A49").Select()
$excel.ActiveWindow.ScrollRow = 49
$excel.ActiveWindow.ScrollColumn = 1
$ws.Range("M59").Select()
